$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Col4a6"
$ws.Cells.Item(2, 3).Value = "Cd93"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.09175499999999999
$ws.Cells.Item(2, 8).Value = 0.275265
$ws.Cells.Item(2, 9).Value = 0.07639436047458779
$ws.Cells.Item(2, 10).Value = 0.07639436047458781
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 281.0920463333333
$ws.Cells.Item(2, 14).Value = 843.2761389999999
$ws.Cells.Item(2, 15).Value = 0.8291026083535286
$ws.Cells.Item(2, 16).Value = 0.8291026083535286
$ws.Cells.Item(2, 17).Value = 25.791600711315
$ws.Cells.Item(2, 18).Value = 232.124406401835
$ws.Cells.Item(2, 19).Value = 0.06333876353298044
$ws.Cells.Item(2, 20).Value = 0.06333876353298046

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Col4a6"
$ws.Cells.Item(3, 3).Value = "Cd93"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.09175499999999999
$ws.Cells.Item(3, 8).Value = 0.275265
$ws.Cells.Item(3, 9).Value = 0.07639436047458779
$ws.Cells.Item(3, 10).Value = 0.07639436047458781
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.449122
$ws.Cells.Item(3, 14).Value = 1.347366
$ws.Cells.Item(3, 15).Value = 0.001324719879221983
$ws.Cells.Item(3, 16).Value = 0.001324719879221983
$ws.Cells.Item(3, 17).Value = 0.04120918911
$ws.Cells.Item(3, 18).Value = 0.37088270199
$ws.Cells.Item(3, 19).Value = 0.0001012011279811366
$ws.Cells.Item(3, 20).Value = 0.0001012011279811366

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Col4a6"
$ws.Cells.Item(4, 3).Value = "Cd93"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.09175499999999999
$ws.Cells.Item(4, 8).Value = 0.275265
$ws.Cells.Item(4, 9).Value = 0.07639436047458779
$ws.Cells.Item(4, 10).Value = 0.07639436047458781
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.452417
$ws.Cells.Item(4, 14).Value = 13.357251
$ws.Cells.Item(4, 15).Value = 0.01313274635953239
$ws.Cells.Item(4, 16).Value = 0.01313274635953239
$ws.Cells.Item(4, 17).Value = 0.4085315218349999
$ws.Cells.Item(4, 18).Value = 3.676783696515
$ws.Cells.Item(4, 19).Value = 0.001003267759411448
$ws.Cells.Item(4, 20).Value = 0.001003267759411448

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Col4a6"
$ws.Cells.Item(5, 3).Value = "Cd93"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.09175499999999999
$ws.Cells.Item(5, 8).Value = 0.275265
$ws.Cells.Item(5, 9).Value = 0.07639436047458779
$ws.Cells.Item(5, 10).Value = 0.07639436047458781
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 53.03808999999999
$ws.Cells.Item(5, 14).Value = 159.11427
$ws.Cells.Item(5, 15).Value = 0.156439925407717
$ws.Cells.Item(5, 16).Value = 0.156439925407717
$ws.Cells.Item(5, 17).Value = 4.866509947949998
$ws.Cells.Item(5, 18).Value = 43.79858953154999
$ws.Cells.Item(5, 19).Value = 0.01195112805421476
$ws.Cells.Item(5, 20).Value = 0.01195112805421476

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Col4a6"
$ws.Cells.Item(6, 3).Value = "Cd93"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.6992243333333333
$ws.Cells.Item(6, 8).Value = 2.097673
$ws.Cells.Item(6, 9).Value = 0.5821676832136669
$ws.Cells.Item(6, 10).Value = 0.5821676832136669
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 281.0920463333333
$ws.Cells.Item(6, 14).Value = 843.2761389999999
$ws.Cells.Item(6, 15).Value = 0.8291026083535286
$ws.Cells.Item(6, 16).Value = 0.8291026083535286
$ws.Cells.Item(6, 17).Value = 196.5463987027274
$ws.Cells.Item(6, 18).Value = 1768.917588324547
$ws.Cells.Item(6, 19).Value = 0.4826767446515819
$ws.Cells.Item(6, 20).Value = 0.4826767446515819

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Col4a6"
$ws.Cells.Item(7, 3).Value = "Cd93"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.6992243333333333
$ws.Cells.Item(7, 8).Value = 2.097673
$ws.Cells.Item(7, 9).Value = 0.5821676832136669
$ws.Cells.Item(7, 10).Value = 0.5821676832136669
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.449122
$ws.Cells.Item(7, 14).Value = 1.347366
$ws.Cells.Item(7, 15).Value = 0.001324719879221983
$ws.Cells.Item(7, 16).Value = 0.001324719879221983
$ws.Cells.Item(7, 17).Value = 0.3140370310353334
$ws.Cells.Item(7, 18).Value = 2.826333279318
$ws.Cells.Item(7, 19).Value = 0.0007712091029937507
$ws.Cells.Item(7, 20).Value = 0.0007712091029937507

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Col4a6"
$ws.Cells.Item(8, 3).Value = "Cd93"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.6992243333333333
$ws.Cells.Item(8, 8).Value = 2.097673
$ws.Cells.Item(8, 9).Value = 0.5821676832136669
$ws.Cells.Item(8, 10).Value = 0.5821676832136669
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 4.452417
$ws.Cells.Item(8, 14).Value = 13.357251
$ws.Cells.Item(8, 15).Value = 0.01313274635953239
$ws.Cells.Item(8, 16).Value = 0.01313274635953239
$ws.Cells.Item(8, 17).Value = 3.113238308547
$ws.Cells.Item(8, 18).Value = 28.019144776923
$ws.Cells.Item(8, 19).Value = 0.007645460522361688
$ws.Cells.Item(8, 20).Value = 0.007645460522361688

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Col4a6"
$ws.Cells.Item(9, 3).Value = "Cd93"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.6992243333333333
$ws.Cells.Item(9, 8).Value = 2.097673
$ws.Cells.Item(9, 9).Value = 0.5821676832136669
$ws.Cells.Item(9, 10).Value = 0.5821676832136669
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 53.03808999999999
$ws.Cells.Item(9, 14).Value = 159.11427
$ws.Cells.Item(9, 15).Value = 0.156439925407717
$ws.Cells.Item(9, 16).Value = 0.156439925407717
$ws.Cells.Item(9, 17).Value = 37.08552312152332
$ws.Cells.Item(9, 18).Value = 333.7697080937099
$ws.Cells.Item(9, 19).Value = 0.09107426893672944
$ws.Cells.Item(9, 20).Value = 0.09107426893672947

$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Col4a6"
$ws.Cells.Item(10, 3).Value = "Cd93"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.410091
$ws.Cells.Item(10, 8).Value = 1.230273
$ws.Cells.Item(10, 9).Value = 0.3414379563117452
$ws.Cells.Item(10, 10).Value = 0.3414379563117453
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 281.0920463333333
$ws.Cells.Item(10, 14).Value = 843.2761389999999
$ws.Cells.Item(10, 15).Value = 0.8291026083535286
$ws.Cells.Item(10, 16).Value = 0.8291026083535286
$ws.Cells.Item(10, 17).Value = 115.273318372883
$ws.Cells.Item(10, 18).Value = 1037.459865355947
$ws.Cells.Item(10, 19).Value = 0.2830871001689661
$ws.Cells.Item(10, 20).Value = 0.2830871001689662

$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Col4a6"
$ws.Cells.Item(11, 3).Value = "Cd93"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.410091
$ws.Cells.Item(11, 8).Value = 1.230273
$ws.Cells.Item(11, 9).Value = 0.3414379563117452
$ws.Cells.Item(11, 10).Value = 0.3414379563117453
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.449122
$ws.Cells.Item(11, 14).Value = 1.347366
$ws.Cells.Item(11, 15).Value = 0.001324719879221983
$ws.Cells.Item(11, 16).Value = 0.001324719879221983
$ws.Cells.Item(11, 17).Value = 0.184180890102
$ws.Cells.Item(11, 18).Value = 1.657628010918
$ws.Cells.Item(11, 19).Value = 0.000452309648247096
$ws.Cells.Item(11, 20).Value = 0.0004523096482470961

$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Col4a6"
$ws.Cells.Item(12, 3).Value = "Cd93"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.410091
$ws.Cells.Item(12, 8).Value = 1.230273
$ws.Cells.Item(12, 9).Value = 0.3414379563117452
$ws.Cells.Item(12, 10).Value = 0.3414379563117453
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 4.452417
$ws.Cells.Item(12, 14).Value = 13.357251
$ws.Cells.Item(12, 15).Value = 0.01313274635953239
$ws.Cells.Item(12, 16).Value = 0.01313274635953239
$ws.Cells.Item(12, 17).Value = 1.825896139947
$ws.Cells.Item(12, 18).Value = 16.433065259523
$ws.Cells.Item(12, 19).Value = 0.004484018077759251
$ws.Cells.Item(12, 20).Value = 0.004484018077759251

$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Col4a6"
$ws.Cells.Item(13, 3).Value = "Cd93"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.410091
$ws.Cells.Item(13, 8).Value = 1.230273
$ws.Cells.Item(13, 9).Value = 0.3414379563117452
$ws.Cells.Item(13, 10).Value = 0.3414379563117453
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 53.03808999999999
$ws.Cells.Item(13, 14).Value = 159.11427
$ws.Cells.Item(13, 15).Value = 0.156439925407717
$ws.Cells.Item(13, 16).Value = 0.156439925407717
$ws.Cells.Item(13, 17).Value = 21.75044336618999
$ws.Cells.Item(13, 18).Value = 195.75399029571
$ws.Cells.Item(13, 19).Value = 0.05341452841677275
$ws.Cells.Item(13, 20).Value = 0.05341452841677277
